$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-18 01:59:28"

# --- 1. Insert two fresh rows right below the header so the two newly
#        scraped items can be populated, pushing all existing data down. ---
$ws.Rows("2:3").Insert()

# New item 1: Java instructor gig (priority score 78)
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(2, 2).Value = "【Java/講師/福岡市内】企業向け新入社員研修のJava講師業務(サブ講師)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5488955"
$ws.Cells.Item(2, 7).Value = 78
$ws.Cells.Item(2, 8).Value = "★Java"

# New item 2: WordPress slider gig (priority score 33)
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(3, 2).Value = "WordPressへ実装するスライダーの制作"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5493927"
$ws.Cells.Item(3, 7).Value = 33
$ws.Cells.Item(3, 8).Value = "○WordPress"

# --- 2. The feed was re-fetched, so every row gets a refreshed "fetched at"
#        timestamp. ---
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 3. Re-sort the whole table by priority score (column G), descending,
#        exactly as the scraper does after merging in new listings. ---
$ws.Range("A1:H20").Sort($ws.Range("G1:G20"), 2)

# --- 4. Wire up hyperlinks for the two rows that are "new" to the
#        worksheet's hyperlink table (rows 19-20 after the sort). ---
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://www.lancers.jp/work/detail/5493449") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://www.lancers.jp/work/detail/5493714") | Out-Null
